$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BB (index 54): "IMPORT_AVERAGE" / EQUIPMENT forecast column ---
# Copy formats from the adjacent BA column so the new header (row 1, bold/
# bordered/centered date style) and new date cell in column A (row 83)
# match the existing staircase-table formatting exactly.
$ws.Range("BA1").Copy() | Out-Null
$ws.Range("BB1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A82").Copy() | Out-Null
$ws.Range("A83").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Header date for the new column (row 1)
$ws.Range("BB1").Value = 45986

# Data values for the new column, rows 2-82 (one quarter further than BA)
$bbValues = @(0.04899821040007168,1.885178963001849,-1.307280175329765,2.75267136845396,-4.471592960313714,-4.942495347902479,-5.434114574907241,4.370618442157621,-0.9953340506219917,5.692238679293155,5.484876523251003,0.1757137213762547,1.626538719806248,2.437334396728659,1.566144859651857,1.556824096124856,0.21169683169569,0.1582599013804469,-0.3253450194449812,-0.15853729477206,0.2054976683197225,0.9595582875050894,1.503966953576466,1.799837015295822,0.6518403676065248,0.6633823054011998,0.923139910942723,0.3869820931359413,2.093916965767463,1.580888475204972,0.3008043112709089,1.199598313222268,0.4774400648527148,1.507463254996111,0.2179116434425623,0.6120689161334525,2.489390679284554,0.5389418434166515,2.40118094791471,1.043009620608657,1.830928398766659,-0.3674870133197601,1.62717758729876,1.117271732844245,0.6525147083449099,1.4,-0.3,-0.3,0.1,-1.138880770453937,-16.88491062648744,9.224715108933083,3.283355339827622,4.432584407022276,2.509693347214139,-0.4381048169788073,4.106981763725997,0.3842995656585515,2.277966437795897,2.507553358214992,-2.399190900254823,-1.115644072253531,0.1689348086957096,-1.814969742946232,-1.586779238813989,-0.2996177924633514,2.757652919539751,-0.08601690538415596,-0.6801011570971838,1.538981993999982,1.68501852020853,0.03331000006224372,0.1716413405801304,0.1716413405801304,0.1716413405801304,0.1716413405801304,0.1716413405801304,0.1716413405801304,0.1716413405801304,0.1716413405801304,0.1716413405801304)
for ($i = 0; $i -lt $bbValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 54).Value = $bbValues[$i]
}

# New final row (83): one more quarter than the previous last row (82)
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = 0.1716413405801304
